# Add two new "Backbones" sheets (Backbones2, Backbones3) at the end of the
# workbook, after the existing "Backbones1" sheet, and fill them with the
# node-id/backbone partition data produced by Main.py.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Backbones2" (new sheet4.xml) -- appended after the last sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Backbones2"

$data2 = @(
    @(8, 60, 88, 97, 100, 1, 57),
    @(81, 53, 10, 44, 42, 38, 17),
    @(51, 31, 41, 16, 18, 59, 26),
    @(87, 27, 91, 7, 65, 62, 83),
    @(68, 95, 48, 19, 36, 85, 14),
    @(2, 56, 12, 92, 61, $null, 28),
    @(32, 35, 15, 54, 30, $null, $null),
    @(33, 4, 11, 94, 72, $null, $null),
    @(73, 49, 69, 9, $null, $null, $null),
    @(82, 58, 64, 86, $null, $null, $null),
    @(6, 34, 96, 79, $null, $null, $null),
    @(90, 78, 3, 71, $null, $null, $null),
    @(50, 93, 84, $null, $null, $null, $null),
    @(22, 67, 63, $null, $null, $null, $null),
    @(76, 46, 43, $null, $null, $null, $null),
    @(20, 23, 21, $null, $null, $null, $null),
    @(75, 55, 37, $null, $null, $null, $null),
    @(70, 74, 89, $null, $null, $null, $null),
    @(47, 80, 29, $null, $null, $null, $null),
    @(45, 5, $null, $null, $null, $null, $null),
    @(13, 24, $null, $null, $null, $null, $null),
    @($null, 25, $null, $null, $null, $null, $null),
    @($null, 40, $null, $null, $null, $null, $null),
    @($null, 66, $null, $null, $null, $null, $null),
    @($null, 98, $null, $null, $null, $null, $null),
    @($null, 77, $null, $null, $null, $null, $null),
    @($null, 99, $null, $null, $null, $null, $null),
    @($null, 52, $null, $null, $null, $null, $null),
    @($null, 39, $null, $null, $null, $null, $null)
)

for ($r = 0; $r -lt $data2.Count; $r++) {
    $rowVals = $data2[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $v = $rowVals[$c]
        if ($v -ne $null) {
            $ws2.Cells.Item($r + 1, $c + 1).Value = $v
        }
    }
}

$hdr2 = $ws2.Range("A1:G1")
$hdr2.Font.Bold = $true
$hdr2.Interior.Color = 13492663

# ---------------------------------------------------------------------------
# Sheet "Backbones3" (new sheet5.xml) -- appended after "Backbones2".
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Backbones3"

$data3 = @(
    @(8, 60, 88, 97, 3, 7, 64),
    @(55, 42, $null, 37, 40, 47, 83),
    @(91, 51, $null, 86, 38, 74, $null),
    @(12, 29, $null, 77, 90, $null, $null),
    @(36, 24, $null, 93, 43, $null, $null),
    @(15, 72, $null, 98, 18, $null, $null),
    @(21, 5, $null, 9, 34, $null, $null),
    @(66, 59, $null, 92, 71, $null, $null),
    @(28, 70, $null, 63, 53, $null, $null),
    @(4, 10, $null, 68, 78, $null, $null),
    @(75, 73, $null, 17, 32, $null, $null),
    @(35, 61, $null, 95, $null, $null, $null),
    @(81, 16, $null, 76, $null, $null, $null),
    @(11, 20, $null, 6, $null, $null, $null),
    @(46, 84, $null, 58, $null, $null, $null),
    @(57, 1, $null, 25, $null, $null, $null),
    @(85, 45, $null, $null, $null, $null, $null),
    @(19, 82, $null, $null, $null, $null, $null),
    @(56, 33, $null, $null, $null, $null, $null),
    @(41, 31, $null, $null, $null, $null, $null),
    @(99, 27, $null, $null, $null, $null, $null),
    @(23, 96, $null, $null, $null, $null, $null),
    @(80, 69, $null, $null, $null, $null, $null),
    @(67, 62, $null, $null, $null, $null, $null),
    @(49, 26, $null, $null, $null, $null, $null),
    @(44, 54, $null, $null, $null, $null, $null),
    @(52, 48, $null, $null, $null, $null, $null),
    @(100, 30, $null, $null, $null, $null, $null),
    @(14, 22, $null, $null, $null, $null, $null),
    @(39, 2, $null, $null, $null, $null, $null),
    @(87, 89, $null, $null, $null, $null, $null),
    @($null, 79, $null, $null, $null, $null, $null),
    @($null, 13, $null, $null, $null, $null, $null),
    @($null, 50, $null, $null, $null, $null, $null),
    @($null, 65, $null, $null, $null, $null, $null),
    @($null, 94, $null, $null, $null, $null, $null)
)

for ($r = 0; $r -lt $data3.Count; $r++) {
    $rowVals = $data3[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $v = $rowVals[$c]
        if ($v -ne $null) {
            $ws3.Cells.Item($r + 1, $c + 1).Value = $v
        }
    }
}

$hdr3 = $ws3.Range("A1:G1")
$hdr3.Font.Bold = $true
$hdr3.Interior.Color = 13492663
